# Auto-generated Excel COM-interop script to update cryptos.xlsx values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells to be stored as text so numeric-looking values
# (e.g. '1.001', '1.000', '113.00') are not reinterpreted as numbers.
$dCells = @('D2','D3','D4','D5','D7','D8','D9','D10','D12','D13','D14','D15','D16','D18','D19','D20','D21','D22','D24','D25','D26','D27','D28','D29','D30','D31','D32','D33','D35','D36','D37','D38','D39','D40','D41','D42','D44','D45','D46','D47','D48','D49','D50','D51')
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range('D2').Value = '29.678.53'
$ws.Range('E2').Value = '  +1.07%  '
$ws.Range('D3').Value = '1.920.69'
$ws.Range('E3').Value = '  -0.44%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.71%  '
$ws.Range('D5').Value = '334.87'
$ws.Range('E5').Value = '  -1.42%  '
$ws.Range('E6').Value = '  -0.64%  '
$ws.Range('D7').Value = '0.4668'
$ws.Range('E7').Value = '  -1.24%  '
$ws.Range('D8').Value = '0.4142'
$ws.Range('E8').Value = '  +1.51%  '
$ws.Range('D9').Value = '48.39'
$ws.Range('E9').Value = '  +1.12%  '
$ws.Range('D10').Value = '0.08065'
$ws.Range('E10').Value = '  -1.21%  '
$ws.Range('E11').Value = '  -0.14%  '
$ws.Range('D12').Value = '22.29'
$ws.Range('E12').Value = '  -0.59%  '
$ws.Range('D13').Value = '1.924.50'
$ws.Range('E13').Value = '  -1.53%  '
$ws.Range('D14').Value = '6.018'
$ws.Range('E14').Value = '  -1.04%  '
$ws.Range('D15').Value = '7.187'
$ws.Range('E15').Value = '  -1.33%  '
$ws.Range('D16').Value = '89.85'
$ws.Range('E16').Value = '  -1.23%  '
$ws.Range('D18').Value = '0.00001038'
$ws.Range('E18').Value = '  -1.17%  '
$ws.Range('D19').Value = '0.06595'
$ws.Range('E19').Value = '  -1.24%  '
$ws.Range('D20').Value = '17.83'
$ws.Range('E20').Value = '  +0.27%  '
$ws.Range('D21').Value = '1.002'
$ws.Range('E21').Value = '  -0.21%  '
$ws.Range('D22').Value = '29.670.98'
$ws.Range('E22').Value = '  +0.98%  '
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('D24').Value = '11.59'
$ws.Range('E24').Value = '  +4.03%  '
$ws.Range('D25').Value = '2.199'
$ws.Range('E25').Value = '  -3.69%  '
$ws.Range('D26').Value = '2.140.34'
$ws.Range('E26').Value = '  -1.57%  '
$ws.Range('D27').Value = '157.85'
$ws.Range('E27').Value = '  -2.08%  '
$ws.Range('D28').Value = '19.94'
$ws.Range('E28').Value = '  -0.24%  '
$ws.Range('D29').Value = '2.153'
$ws.Range('E29').Value = '  +0.50%  '
$ws.Range('D30').Value = '5.727'
$ws.Range('E30').Value = '  +2.80%  '
$ws.Range('D31').Value = '117.98'
$ws.Range('E31').Value = '  -3.07%  '
$ws.Range('D32').Value = '1.048'
$ws.Range('E32').Value = '  +5.46%  '
$ws.Range('D33').Value = '0.09451'
$ws.Range('E33').Value = '  -1.36%  '
$ws.Range('E34').Value = '  -1.21%  '
$ws.Range('D35').Value = '5.441'
$ws.Range('E35').Value = '  +0.49%  '
$ws.Range('D36').Value = '3.526'
$ws.Range('E36').Value = '  -3.75%  '
$ws.Range('D37').Value = '0.06147'
$ws.Range('E37').Value = '  -1.03%  '
$ws.Range('D38').Value = '0.02268'
$ws.Range('E38').Value = '  -0.94%  '
$ws.Range('D39').Value = '8.468'
$ws.Range('E39').Value = '  +0.30%  '
$ws.Range('D40').Value = '1.179'
$ws.Range('E40').Value = '  -0.17%  '
$ws.Range('D41').Value = '0.5916'
$ws.Range('E41').Value = '  -1.61%  '
$ws.Range('D42').Value = '1.000'
$ws.Range('E42').Value = '  -0.63%  '
$ws.Range('E43').Value = '  -1.81%  '
$ws.Range('D44').Value = '10.26'
$ws.Range('E44').Value = '  -2.63%  '
$ws.Range('D45').Value = '2.337'
$ws.Range('E45').Value = '  +0.19%  '
$ws.Range('D46').Value = '1.240'
$ws.Range('E46').Value = '  -1.10%  '
$ws.Range('D47').Value = '0.07519'
$ws.Range('E47').Value = '  +2.56%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '12.28'
$ws.Range('E48').Value = '  -0.93%  '
$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').Value = '0.5588'
$ws.Range('E49').Value = '  -1.21%  '
$ws.Range('D50').Value = '1.943'
$ws.Range('E50').Value = '  -1.32%  '
$ws.Range('D51').Value = '113.00'
$ws.Range('E51').Value = '  +0.52%  '
